# Shrink the font size of several "Rounded Rectangle" callout labels from
# 14pt to 13.5pt, and tidy up the "Kidney crystallopathy ..." box so the
# trailing star-rating text shares a single run with the preceding text
# (both end up at 13.5pt).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "Neural tube closure defects  ★★★★☆" (single run) : 14 -> 13.5
$s.Shapes.Item(14).TextFrame.TextRange.Font.Size = 13.5

# "Acute inflammation" / "★★★★★" (two paragraphs, one run each) : 14 -> 13.5
$s.Shapes.Item(15).TextFrame.TextRange.Font.Size = 13.5

# "Kidney crystallopathy and tubular necrosis ★☆☆☆☆" : the last two runs
# (" and tubular necrosis " at 13.5 and "★☆☆☆☆" at 14) collapse into a
# single run at 13.5pt. Replacing the text of the combined character span
# with itself (same text) merges the underlying runs and adopts the
# formatting of the first run in that span.
$tr40 = $s.Shapes.Item(16).TextFrame.TextRange
$tail = $tr40.Characters(22, 27)
$tail.Text = " and tubular necrosis ★☆☆☆☆"

# "Cognitive function defects  ★☆☆☆☆" (single run) : 14 -> 13.5
$s.Shapes.Item(17).TextFrame.TextRange.Font.Size = 13.5

# "Immune-mediated diseases  ★☆☆☆☆" (single run) : 14 -> 13.5
$s.Shapes.Item(18).TextFrame.TextRange.Font.Size = 13.5

# "Chronic inflammation" / "★★★☆☆" (two paragraphs, one run each) : 14 -> 13.5
$s.Shapes.Item(19).TextFrame.TextRange.Font.Size = 13.5
